$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Month) was converted from month names (July, August, ...)
# to 2-digit month-number codes (text, formatted with the Text number
# format so leading zeros like "07" are preserved). Oct/Nov/Dec (10/11/12)
# don't need a leading zero, so several of those rows were entered as plain
# numbers instead - still with the Text number format applied afterwards.
$monthData = @{
  2 = @{ code = "07"; isstr = $true }
  3 = @{ code = "08"; isstr = $true }
  4 = @{ code = "09"; isstr = $true }
  5 = @{ code = "10"; isstr = $false }
  6 = @{ code = "11"; isstr = $false }
  7 = @{ code = "12"; isstr = $false }
  8 = @{ code = "01"; isstr = $true }
  9 = @{ code = "02"; isstr = $true }
  10 = @{ code = "03"; isstr = $true }
  11 = @{ code = "04"; isstr = $true }
  12 = @{ code = "05"; isstr = $true }
  13 = @{ code = "06"; isstr = $true }
  14 = @{ code = "07"; isstr = $true }
  15 = @{ code = "08"; isstr = $true }
  16 = @{ code = "09"; isstr = $true }
  17 = @{ code = "10"; isstr = $false }
  18 = @{ code = "11"; isstr = $false }
  19 = @{ code = "12"; isstr = $false }
  20 = @{ code = "01"; isstr = $true }
  21 = @{ code = "02"; isstr = $true }
  22 = @{ code = "03"; isstr = $true }
  23 = @{ code = "04"; isstr = $true }
  24 = @{ code = "05"; isstr = $true }
  25 = @{ code = "06"; isstr = $true }
  26 = @{ code = "07"; isstr = $true }
  27 = @{ code = "08"; isstr = $true }
  28 = @{ code = "09"; isstr = $true }
  29 = @{ code = "10"; isstr = $false }
  30 = @{ code = "11"; isstr = $false }
  31 = @{ code = "12"; isstr = $true }
  32 = @{ code = "01"; isstr = $true }
  33 = @{ code = "02"; isstr = $true }
  34 = @{ code = "03"; isstr = $true }
  35 = @{ code = "04"; isstr = $true }
  36 = @{ code = "05"; isstr = $true }
  37 = @{ code = "06"; isstr = $true }
  38 = @{ code = "07"; isstr = $true }
  39 = @{ code = "08"; isstr = $true }
  40 = @{ code = "09"; isstr = $true }
  41 = @{ code = "10"; isstr = $true }
  42 = @{ code = "11"; isstr = $true }
  43 = @{ code = "12"; isstr = $false }
  44 = @{ code = "01"; isstr = $true }
  45 = @{ code = "02"; isstr = $true }
  46 = @{ code = "03"; isstr = $true }
  47 = @{ code = "04"; isstr = $true }
  48 = @{ code = "10"; isstr = $true }
  49 = @{ code = "11"; isstr = $false }
  50 = @{ code = "12"; isstr = $false }
  51 = @{ code = "01"; isstr = $true }
  52 = @{ code = "02"; isstr = $true }
  53 = @{ code = "03"; isstr = $true }
  54 = @{ code = "04"; isstr = $true }
  55 = @{ code = "05"; isstr = $true }
  56 = @{ code = "06"; isstr = $true }
  57 = @{ code = "07"; isstr = $true }
  58 = @{ code = "08"; isstr = $true }
  59 = @{ code = "09"; isstr = $true }
  60 = @{ code = "10"; isstr = $true }
  61 = @{ code = "11"; isstr = $false }
  62 = @{ code = "12"; isstr = $false }
  63 = @{ code = "01"; isstr = $true }
  64 = @{ code = "02"; isstr = $true }
  65 = @{ code = "03"; isstr = $true }
  66 = @{ code = "04"; isstr = $true }
  67 = @{ code = "05"; isstr = $true }
  68 = @{ code = "06"; isstr = $true }
  69 = @{ code = "07"; isstr = $true }
  70 = @{ code = "08"; isstr = $true }
  71 = @{ code = "09"; isstr = $true }
  72 = @{ code = "10"; isstr = $true }
  73 = @{ code = "11"; isstr = $false }
  74 = @{ code = "12"; isstr = $false }
  75 = @{ code = "01"; isstr = $true }
  76 = @{ code = "02"; isstr = $true }
  77 = @{ code = "03"; isstr = $true }
  78 = @{ code = "04"; isstr = $true }
  79 = @{ code = "05"; isstr = $true }
  80 = @{ code = "06"; isstr = $true }
  81 = @{ code = "07"; isstr = $true }
  82 = @{ code = "08"; isstr = $true }
  83 = @{ code = "07"; isstr = $true }
  84 = @{ code = "08"; isstr = $true }
  85 = @{ code = "09"; isstr = $true }
  86 = @{ code = "10"; isstr = $true }
  87 = @{ code = "11"; isstr = $false }
  88 = @{ code = "12"; isstr = $false }
  89 = @{ code = "01"; isstr = $true }
  90 = @{ code = "02"; isstr = $true }
  91 = @{ code = "03"; isstr = $true }
  92 = @{ code = "04"; isstr = $true }
  93 = @{ code = "05"; isstr = $true }
  94 = @{ code = "06"; isstr = $true }
  95 = @{ code = "07"; isstr = $true }
  96 = @{ code = "08"; isstr = $true }
  97 = @{ code = "09"; isstr = $true }
  98 = @{ code = "10"; isstr = $false }
  99 = @{ code = "11"; isstr = $false }
  100 = @{ code = "12"; isstr = $false }
  101 = @{ code = "01"; isstr = $true }
  102 = @{ code = "02"; isstr = $true }
  103 = @{ code = "03"; isstr = $true }
  104 = @{ code = "04"; isstr = $true }
  105 = @{ code = "05"; isstr = $true }
  106 = @{ code = "06"; isstr = $true }
  107 = @{ code = "07"; isstr = $true }
  108 = @{ code = "08"; isstr = $true }
  109 = @{ code = "09"; isstr = $true }
  110 = @{ code = "10"; isstr = $true }
  111 = @{ code = "11"; isstr = $true }
  112 = @{ code = "12"; isstr = $false }
  113 = @{ code = "01"; isstr = $true }
  114 = @{ code = "02"; isstr = $true }
  115 = @{ code = "03"; isstr = $true }
  116 = @{ code = "04"; isstr = $true }
  117 = @{ code = "05"; isstr = $true }
  118 = @{ code = "06"; isstr = $true }
  119 = @{ code = "07"; isstr = $true }
  120 = @{ code = "08"; isstr = $true }
  121 = @{ code = "11"; isstr = $false }
  122 = @{ code = "12"; isstr = $false }
  123 = @{ code = "01"; isstr = $true }
  124 = @{ code = "02"; isstr = $true }
  125 = @{ code = "03"; isstr = $true }
  126 = @{ code = "04"; isstr = $true }
  127 = @{ code = "05"; isstr = $true }
  128 = @{ code = "06"; isstr = $true }
  129 = @{ code = "07"; isstr = $true }
  130 = @{ code = "08"; isstr = $true }
  131 = @{ code = "09"; isstr = $true }
  132 = @{ code = "10"; isstr = $true }
  133 = @{ code = "11"; isstr = $false }
  134 = @{ code = "12"; isstr = $false }
  135 = @{ code = "01"; isstr = $true }
  136 = @{ code = "02"; isstr = $true }
  137 = @{ code = "03"; isstr = $true }
  138 = @{ code = "04"; isstr = $true }
  139 = @{ code = "05"; isstr = $true }
  140 = @{ code = "06"; isstr = $true }
  141 = @{ code = "07"; isstr = $true }
  142 = @{ code = "08"; isstr = $true }
  143 = @{ code = "09"; isstr = $true }
  144 = @{ code = "10"; isstr = $false }
  145 = @{ code = "11"; isstr = $false }
  146 = @{ code = "12"; isstr = $false }
  147 = @{ code = "01"; isstr = $true }
  148 = @{ code = "02"; isstr = $true }
  149 = @{ code = "03"; isstr = $true }
  150 = @{ code = "04"; isstr = $true }
  151 = @{ code = "05"; isstr = $true }
  152 = @{ code = "06"; isstr = $true }
  153 = @{ code = "07"; isstr = $true }
}

foreach ($row in $monthData.Keys) {
  $info = $monthData[$row]
  $cell = $ws.Range("A" + $row)
  if ($info.isstr) {
    # Leading-zero codes ("01".."09") must be entered as text, so the Text
    # number format is set BEFORE assigning the value (otherwise the
    # numeric-looking string would be coerced back into a number).
    $cell.NumberFormat = "@"
    $cell.Value = $info.code
  } else {
    # No leading zero needed (10/11/12) - assign as a real number first,
    # then switch the cell to the Text format (keeps it stored as a number).
    $cell.Value = [int]$info.code
    $cell.NumberFormat = "@"
  }
}

# Scroll/selection state left by the editor
$ws.Activate()
[void]$ws.Range("H13").Select()

# Page setup: paper size + orientation changed
$pageSetup = $ws.PageSetup
$pageSetup.PaperSize = 150
$pageSetup.Orientation = 1

